# Update "想去人数" (interest count) and "最低票价" (min ticket price) figures
# on both the "展览" and "全部类型" worksheets, mirroring a data refresh.

$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

# Cell -> new value map (identical for both sheets)
$updates = @{
    "G2"  = 65
    "F10" = 44
    "F11" = 18
    "F12" = 556
    "F14" = 288
    "F16" = 348
    "F22" = 868
    "F23" = 1381
    "F24" = 290
    "F25" = 315
    "F28" = 152
    "F29" = 35
    "F30" = 84
    "F31" = 207
    "F33" = 261
    "F34" = 1600
    "F37" = 153
    "F40" = 3487
    "F41" = 412
    "F42" = 182
    "F43" = 888
}

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($cellRef in $updates.Keys) {
        $ws.Range($cellRef).Value = $updates[$cellRef]
    }
}
